$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 61.06878433333333
$ws.Range("H2").Value = 183.206353
$ws.Range("I2").Value = 0.4308066250287063
$ws.Range("J2").Value = 0.4308066250287063
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.222961333333333
$ws.Range("N2").Value = 3.668884
$ws.Range("O2").Value = 0.02082890820948692
$ws.Range("P2").Value = 0.02082890820948692
$ws.Range("Q2").Value = 74.68476191333912
$ws.Range("R2").Value = 672.162857220052
$ws.Range("S2").Value = 0.008973231648761773
$ws.Range("T2").Value = 0.008973231648761777

# Row 3
$ws.Range("G3").Value = 61.06878433333333
$ws.Range("H3").Value = 183.206353
$ws.Range("I3").Value = 0.4308066250287063
$ws.Range("J3").Value = 0.4308066250287063
$ws.Range("O3").Value = 0.1691870972318839
$ws.Range("P3").Value = 0.169187097231884
$ws.Range("Q3").Value = 606.642362071433
$ws.Range("R3").Value = 5459.781258642898
$ws.Range("S3").Value = 0.0728869223568715
$ws.Range("T3").Value = 0.07288692235687153

# Row 4
$ws.Range("G4").Value = 61.06878433333333
$ws.Range("H4").Value = 183.206353
$ws.Range("I4").Value = 0.4308066250287063
$ws.Range("J4").Value = 0.4308066250287063
$ws.Range("M4").Value = 14.516908
$ws.Range("N4").Value = 43.550724
$ws.Range("O4").Value = 0.247245220250272
$ws.Range("P4").Value = 0.2472452202502721
$ws.Range("Q4").Value = 886.5299238388413
$ws.Range("R4").Value = 7978.769314549573
$ws.Range("S4").Value = 0.1065148788904988
$ws.Range("T4").Value = 0.1065148788904989

# Row 5
$ws.Range("G5").Value = 61.06878433333333
$ws.Range("H5").Value = 183.206353
$ws.Range("I5").Value = 0.4308066250287063
$ws.Range("J5").Value = 0.4308066250287063
$ws.Range("M5").Value = 12.24131666666667
$ws.Range("N5").Value = 36.72395
$ws.Range("O5").Value = 0.2084884078209579
$ws.Range("P5").Value = 0.2084884078209579
$ws.Range("Q5").Value = 747.5623274727056
$ws.Range("R5").Value = 6728.060947254351
$ws.Range("S5").Value = 0.0898181873309554
$ws.Range("T5").Value = 0.08981818733095541

# Row 6
$ws.Range("G6").Value = 61.06878433333333
$ws.Range("H6").Value = 183.206353
$ws.Range("I6").Value = 0.4308066250287063
$ws.Range("J6").Value = 0.4308066250287063
$ws.Range("M6").Value = 20.799674
$ws.Range("N6").Value = 62.399022
$ws.Range("O6").Value = 0.3542503664873991
$ws.Range("P6").Value = 0.3542503664873992
$ws.Range("Q6").Value = 1270.210805709641
$ws.Range("R6").Value = 11431.89725138677
$ws.Range("S6").Value = 0.1526134048016187
$ws.Range("T6").Value = 0.1526134048016188

# Row 7
$ws.Range("I7").Value = 0.02554841368886107
$ws.Range("J7").Value = 0.02554841368886107
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 1.222961333333333
$ws.Range("N7").Value = 3.668884
$ws.Range("O7").Value = 0.02082890820948692
$ws.Range("P7").Value = 0.02082890820948692
$ws.Range("Q7").Value = 4.429080433684001
$ws.Range("R7").Value = 39.861723903156
$ws.Range("S7").Value = 0.0005321455636232864
$ws.Range("T7").Value = 0.0005321455636232865

# Row 8
$ws.Range("I8").Value = 0.02554841368886107
$ws.Range("J8").Value = 0.02554841368886107
$ws.Range("O8").Value = 0.1691870972318839
$ws.Range("P8").Value = 0.169187097231884
$ws.Range("S8").Value = 0.004322461950897733
$ws.Range("T8").Value = 0.004322461950897734

# Row 9
$ws.Range("I9").Value = 0.02554841368886107
$ws.Range("J9").Value = 0.02554841368886107
$ws.Range("M9").Value = 14.516908
$ws.Range("N9").Value = 43.550724
$ws.Range("O9").Value = 0.247245220250272
$ws.Range("P9").Value = 0.2472452202502721
$ws.Range("Q9").Value = 52.57447756352401
$ws.Range("R9").Value = 473.170298071716
$ws.Range("S9").Value = 0.006316723169547521
$ws.Range("T9").Value = 0.006316723169547522

# Row 10
$ws.Range("I10").Value = 0.02554841368886107
$ws.Range("J10").Value = 0.02554841368886107
$ws.Range("M10").Value = 12.24131666666667
$ws.Range("N10").Value = 36.72395
$ws.Range("O10").Value = 0.2084884078209579
$ws.Range("P10").Value = 0.2084884078209579
$ws.Range("Q10").Value = 44.33318916395001
$ws.Range("R10").Value = 398.99870247555
$ws.Range("S10").Value = 0.005326548092341811
$ws.Range("T10").Value = 0.005326548092341811

# Row 11
$ws.Range("I11").Value = 0.02554841368886107
$ws.Range("J11").Value = 0.02554841368886107
$ws.Range("M11").Value = 20.799674
$ws.Range("N11").Value = 62.399022
$ws.Range("O11").Value = 0.3542503664873991
$ws.Range("P11").Value = 0.3542503664873992
$ws.Range("Q11").Value = 75.32816175742201
$ws.Range("R11").Value = 677.9534558167981
$ws.Range("S11").Value = 0.00905053491245072
$ws.Range("T11").Value = 0.009050534912450722

# Row 12
$ws.Range("G12").Value = 36.843258
$ws.Range("H12").Value = 110.529774
$ws.Range("I12").Value = 0.2599088848306786
$ws.Range("J12").Value = 0.2599088848306786
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 1.222961333333333
$ws.Range("N12").Value = 3.668884
$ws.Range("O12").Value = 0.02082890820948692
$ws.Range("P12").Value = 0.02082890820948692
$ws.Range("Q12").Value = 45.057879928024
$ws.Range("R12").Value = 405.520919352216
$ws.Range("S12").Value = 0.005413618304968312
$ws.Range("T12").Value = 0.005413618304968314

# Row 13
$ws.Range("G13").Value = 36.843258
$ws.Range("H13").Value = 110.529774
$ws.Range("I13").Value = 0.2599088848306786
$ws.Range("J13").Value = 0.2599088848306786
$ws.Range("O13").Value = 0.1691870972318839
$ws.Range("P13").Value = 0.169187097231884
$ws.Range("Q13").Value = 365.9919106548759
$ws.Range("R13").Value = 3293.927195893884
$ws.Range("S13").Value = 0.04397322976927855
$ws.Range("T13").Value = 0.04397322976927857

# Row 14
$ws.Range("G14").Value = 36.843258
$ws.Range("H14").Value = 110.529774
$ws.Range("I14").Value = 0.2599088848306786
$ws.Range("J14").Value = 0.2599088848306786
$ws.Range("M14").Value = 14.516908
$ws.Range("N14").Value = 43.550724
$ws.Range("O14").Value = 0.247245220250272
$ws.Range("P14").Value = 0.2472452202502721
$ws.Range("Q14").Value = 534.850186806264
$ws.Range("R14").Value = 4813.651681256377
$ws.Range("S14").Value = 0.06426122947496372
$ws.Range("T14").Value = 0.06426122947496374

# Row 15
$ws.Range("G15").Value = 36.843258
$ws.Range("H15").Value = 110.529774
$ws.Range("I15").Value = 0.2599088848306786
$ws.Range("J15").Value = 0.2599088848306786
$ws.Range("M15").Value = 12.24131666666667
$ws.Range("N15").Value = 36.72395
$ws.Range("O15").Value = 0.2084884078209579
$ws.Range("P15").Value = 0.2084884078209579
$ws.Range("Q15").Value = 451.0099882097001
$ws.Range("R15").Value = 4059.0898938873
$ws.Range("S15").Value = 0.0541879895768689
$ws.Range("T15").Value = 0.05418798957686891

# Row 16
$ws.Range("G16").Value = 36.843258
$ws.Range("H16").Value = 110.529774
$ws.Range("I16").Value = 0.2599088848306786
$ws.Range("J16").Value = 0.2599088848306786
$ws.Range("M16").Value = 20.799674
$ws.Range("N16").Value = 62.399022
$ws.Range("O16").Value = 0.3542503664873991
$ws.Range("P16").Value = 0.3542503664873992
$ws.Range("Q16").Value = 766.327755497892
$ws.Range("R16").Value = 6896.949799481028
$ws.Range("S16").Value = 0.09207281770459912
$ws.Range("T16").Value = 0.09207281770459914

# Row 17
$ws.Range("G17").Value = 2.119603
$ws.Range("H17").Value = 6.358808999999999
$ws.Range("I17").Value = 0.01495263127961596
$ws.Range("J17").Value = 0.01495263127961596
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 1.222961333333333
$ws.Range("N17").Value = 3.668884
$ws.Range("O17").Value = 0.02082890820948692
$ws.Range("P17").Value = 0.02082890820948692
$ws.Range("Q17").Value = 2.592192511017333
$ws.Range("R17").Value = 23.329732599156
$ws.Range("S17").Value = 0.0003114469844134237
$ws.Range("T17").Value = 0.0003114469844134238

# Row 18
$ws.Range("G18").Value = 2.119603
$ws.Range("H18").Value = 6.358808999999999
$ws.Range("I18").Value = 0.01495263127961596
$ws.Range("J18").Value = 0.01495263127961596
$ws.Range("O18").Value = 0.1691870972318839
$ws.Range("P18").Value = 0.169187097231884
$ws.Range("Q18").Value = 21.05561760579933
$ws.Range("R18").Value = 189.500558452194
$ws.Range("S18").Value = 0.002529792282176894
$ws.Range("T18").Value = 0.002529792282176895

# Row 19
$ws.Range("G19").Value = 2.119603
$ws.Range("H19").Value = 6.358808999999999
$ws.Range("I19").Value = 0.01495263127961596
$ws.Range("J19").Value = 0.01495263127961596
$ws.Range("M19").Value = 14.516908
$ws.Range("N19").Value = 43.550724
$ws.Range("O19").Value = 0.247245220250272
$ws.Range("P19").Value = 0.2472452202502721
$ws.Range("Q19").Value = 30.770081747524
$ws.Range("R19").Value = 276.930735727716
$ws.Range("S19").Value = 0.003696966614049754
$ws.Range("T19").Value = 0.003696966614049755

# Row 20
$ws.Range("G20").Value = 2.119603
$ws.Range("H20").Value = 6.358808999999999
$ws.Range("I20").Value = 0.01495263127961596
$ws.Range("J20").Value = 0.01495263127961596
$ws.Range("M20").Value = 12.24131666666667
$ws.Range("N20").Value = 36.72395
$ws.Range("O20").Value = 0.2084884078209579
$ws.Range("P20").Value = 0.2084884078209579
$ws.Range("Q20").Value = 25.94673153061667
$ws.Range("R20").Value = 233.52058377555
$ws.Range("S20").Value = 0.003117450288220983
$ws.Range("T20").Value = 0.003117450288220983

# Row 21
$ws.Range("G21").Value = 2.119603
$ws.Range("H21").Value = 6.358808999999999
$ws.Range("I21").Value = 0.01495263127961596
$ws.Range("J21").Value = 0.01495263127961596
$ws.Range("M21").Value = 20.799674
$ws.Range("N21").Value = 62.399022
$ws.Range("O21").Value = 0.3542503664873991
$ws.Range("P21").Value = 0.3542503664873992
$ws.Range("Q21").Value = 44.087051409422
$ws.Range("R21").Value = 396.7834626847979
$ws.Range("S21").Value = 0.0052969751107549
$ws.Range("T21").Value = 0.005296975110754902

# Row 22
$ws.Range("G22").Value = 38.101267
$ws.Range("H22").Value = 114.303801
$ws.Range("I22").Value = 0.268783445172138
$ws.Range("J22").Value = 0.268783445172138
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 1.222961333333333
$ws.Range("N22").Value = 3.668884
$ws.Range("O22").Value = 0.02082890820948692
$ws.Range("P22").Value = 0.02082890820948692
$ws.Range("Q22").Value = 46.59637629200933
$ws.Range("R22").Value = 419.367386628084
$ws.Range("S22").Value = 0.005598465707720122
$ws.Range("T22").Value = 0.005598465707720123

# Row 23
$ws.Range("G23").Value = 38.101267
$ws.Range("H23").Value = 114.303801
$ws.Range("I23").Value = 0.268783445172138
$ws.Range("J23").Value = 0.268783445172138
$ws.Range("O23").Value = 0.1691870972318839
$ws.Range("P23").Value = 0.169187097231884
$ws.Range("Q23").Value = 378.4886642680073
$ws.Range("R23").Value = 3406.397978412066
$ws.Range("S23").Value = 0.04547469087265926
$ws.Range("T23").Value = 0.04547469087265927

# Row 24
$ws.Range("G24").Value = 38.101267
$ws.Range("H24").Value = 114.303801
$ws.Range("I24").Value = 0.268783445172138
$ws.Range("J24").Value = 0.268783445172138
$ws.Range("M24").Value = 14.516908
$ws.Range("N24").Value = 43.550724
$ws.Range("O24").Value = 0.247245220250272
$ws.Range("P24").Value = 0.2472452202502721
$ws.Range("Q24").Value = 553.112587722436
$ws.Range("R24").Value = 4978.013289501924
$ws.Range("S24").Value = 0.06645542210121218
$ws.Range("T24").Value = 0.06645542210121218

# Row 25
$ws.Range("G25").Value = 38.101267
$ws.Range("H25").Value = 114.303801
$ws.Range("I25").Value = 0.268783445172138
$ws.Range("J25").Value = 0.268783445172138
$ws.Range("M25").Value = 12.24131666666667
$ws.Range("N25").Value = 36.72395
$ws.Range("O25").Value = 0.2084884078209579
$ws.Range("P25").Value = 0.2084884078209579
$ws.Range("Q25").Value = 466.4096747482167
$ws.Range("R25").Value = 4197.68707273395
$ws.Range("S25").Value = 0.05603823253257079
$ws.Range("T25").Value = 0.05603823253257078

# Row 26
$ws.Range("G26").Value = 38.101267
$ws.Range("H26").Value = 114.303801
$ws.Range("I26").Value = 0.268783445172138
$ws.Range("J26").Value = 0.268783445172138
$ws.Range("M26").Value = 20.799674
$ws.Range("N26").Value = 62.399022
$ws.Range("O26").Value = 0.3542503664873991
$ws.Range("P26").Value = 0.3542503664873992
$ws.Range("Q26").Value = 792.493932586958
$ws.Range("R26").Value = 7132.445393282622
$ws.Range("S26").Value = 0.09521663395797564
$ws.Range("T26").Value = 0.09521663395797565
